$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row (17) of trade data to the sheet, mirroring the existing rows.
$ws.Range("A17").Value = 10227.57
$ws.Range("B17").Value = 9939.33
$ws.Range("C17").Value = 294.14
$ws.Range("D17").Value = 302.66000000000003
$ws.Range("E17").Value = $false
$ws.Range("F17").Value = 2.9
$ws.Range("G17").Value = 42626.544398148151
$ws.Range("G17").NumberFormat = "m/d/yy h:mm"
$ws.Range("H17").Value = $true
